$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textBlob = 'It looks like your message just says "text". Did you have a specific question or topic you''d like to discuss regarding text or anything else? Feel free to let me know how I can assist you!It looks like your message just says "text". Did you have a specifi'
$ingredient = 'Knowledge of pharmaceutical equipments, Medicine effects, Behavior, Functionality, etc.'
$specIO = 'Double-check where the date string is coming from. Ensure that any JavaScript or frontend frameworks formatting dates are consistent where the date string is coming from. Ensure that any JavaScript or frontend frameworks formatting dates are consistent'
$specDetails = 'Observation of bioplant on the presence of WHO.'

foreach ($r in 2..4) {
    $ws.Cells.Item($r, 1).Value = "form_23.pdf"
    $ws.Cells.Item($r, 2).Value = $textBlob
    $ws.Cells.Item($r, 3).Value = $textBlob
    $ws.Cells.Item($r, 4).Value = $textBlob
    $ws.Cells.Item($r, 5).Value = "Volume 1"
    $ws.Cells.Item($r, 6).Value = $ingredient
    $ws.Cells.Item($r, 7).Value = $specIO
    $ws.Cells.Item($r, 8).Value = $specDetails
    $ws.Cells.Item($r, 9).Value = $textBlob
    $ws.Cells.Item($r, 10).Value = $textBlob
}

# Column K (Conclusion) is cleared for row 2 to match the now-blank state
# already present for rows 3-4.
$ws.Cells.Item(2, 11).Value = ""
